$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leaderboard data: Name, Wins, Losses, Ratio
$data = @(
    @("Phillip", 1, 6, 0.167),
    @("Jack",    6, 5, 1.2),
    @("Lance ",  3, 3, 1),
    @("Lance",   4, 1, 4),
    @("Connor",  3, 1, 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
